# ---------------------------------------------------------------------------
# multival_tabular_tests_generated.xlsx regeneration
#
# Starting point (before.xlsx):
#   Sheet1 "Person"  : first_name | last_name | hobbies
#   Sheet2 "Person1" : first_name | last_name | hobbies
#
# Target (per diff):
#   Database      : person_set | pet_set | org_set
#   Organization  : org_name | pet_names
#   Person        : first_name | last_name | hobbies | pet_names  (+ list validation on C)
#   Pet           : pet_name | species
#   Database1     : person_set | pet_set | org_set
#   Organization1 : org_name | pet_names
#   Person1       : first_name | last_name | hobbies | pet_names
#   Pet1          : pet_name | species
#   active tab (0-based) = 7 -> "Pet1"
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Create the six brand-new sheets first, appended at the end of the workbook
# in their target left-to-right relative order, then reorder with Move().
# (Re-fetching worksheets by Item("Name")/Item(Count) right before each call
# matters: this COM layer resolves a previously-captured worksheet handle by
# its *position at the time of the call*, so a stale variable silently
# re-targets after any Add()/Move() that shifts indices.)
$newNames = @("Database", "Organization", "Pet", "Database1", "Organization1", "Pet1")
foreach ($n in $newNames) {
    $s = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
    $s.Name = $n
    # Match the page margins used throughout the rest of the workbook
    # (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 points).
    $s.PageSetup.LeftMargin = 54
    $s.PageSetup.RightMargin = 54
    $s.PageSetup.TopMargin = 72
    $s.PageSetup.BottomMargin = 72
    $s.PageSetup.HeaderMargin = 36
    $s.PageSetup.FooterMargin = 36
}

# Move the two pre-existing sheets into their final slots so the tab order
# becomes: Database, Organization, Person, Pet, Database1, Organization1,
# Person1, Pet1
$wb.Worksheets.Item("Person").Move($wb.Worksheets.Item("Pet"))
$wb.Worksheets.Item("Person1").Move($wb.Worksheets.Item("Pet1"))

# ---------------------------------------------------------------------------
# Sheet contents
# ---------------------------------------------------------------------------

# Database
$ws = $wb.Worksheets.Item("Database")
$ws.Range("A1").Value = "person_set"
$ws.Range("B1").Value = "pet_set"
$ws.Range("C1").Value = "org_set"

# Organization
$ws = $wb.Worksheets.Item("Organization")
$ws.Range("A1").Value = "org_name"
$ws.Range("B1").Value = "pet_names"

# Person (existing sheet -> add 4th column + list validation on hobbies col)
$ws = $wb.Worksheets.Item("Person")
$ws.Range("A1").Value = "first_name"
$ws.Range("B1").Value = "last_name"
$ws.Range("C1").Value = "hobbies"
$ws.Range("D1").Value = "pet_names"
$ws.Range("C2:C1048576").Validation.Add(3, 1, 1, '"tennis,cooking,sewing,fishing"')

# Pet
$ws = $wb.Worksheets.Item("Pet")
$ws.Range("A1").Value = "pet_name"
$ws.Range("B1").Value = "species"

# Database1
$ws = $wb.Worksheets.Item("Database1")
$ws.Range("A1").Value = "person_set"
$ws.Range("B1").Value = "pet_set"
$ws.Range("C1").Value = "org_set"

# Organization1
$ws = $wb.Worksheets.Item("Organization1")
$ws.Range("A1").Value = "org_name"
$ws.Range("B1").Value = "pet_names"

# Person1 (existing sheet -> add 4th column, no validation)
$ws = $wb.Worksheets.Item("Person1")
$ws.Range("A1").Value = "first_name"
$ws.Range("B1").Value = "last_name"
$ws.Range("C1").Value = "hobbies"
$ws.Range("D1").Value = "pet_names"

# Pet1
$ws = $wb.Worksheets.Item("Pet1")
$ws.Range("A1").Value = "pet_name"
$ws.Range("B1").Value = "species"

# ---------------------------------------------------------------------------
# Active tab -> Pet1 (8th sheet, 0-based index 7)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Pet1").Activate()
